# Automatic update of files.
# The underlying observation records for rows 25-30 were re-ordered/re-matched;
# apply the resulting per-row field changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 (now holds the "Tretåig hackspett" / Picoides tridactylus record)
$ws.Range("A25").Value = 111936793
$ws.Range("B25").Value = 56398
$ws.Range("E25").Value = 100109
$ws.Range("F25").Value = "Tretåig hackspett"
$ws.Range("G25").Value = "Picoides tridactylus"
$ws.Range("K25").Value = ""
$ws.Range("L25").Value = ""
$ws.Range("M25").Value = ""
$ws.Range("N25").Value = ""
$ws.Range("Q25").Value = 451088.7179698629
$ws.Range("R25").Value = 7087232.506422138
$ws.Range("AC25").Value = "ringhack äldre"

# Row 26 (now holds the "Doftticka" / Haploporus odorus record)
$ws.Range("A26").Value = 111936894
$ws.Range("B26").Value = 89965
$ws.Range("D26").Value = "VU"
$ws.Range("E26").Value = 760
$ws.Range("F26").Value = "Doftticka"
$ws.Range("G26").Value = "Haploporus odorus"
$ws.Range("H26").Value = "(Sommerf.) Bondartsev & Singer"
$ws.Range("Q26").Value = 451168.6101546783
$ws.Range("R26").Value = 7086616.526546557

# Row 27 (now holds the "Garnlav" / Alectoria sarmentosa record)
$ws.Range("A27").Value = 111936892
$ws.Range("B27").Value = 77515
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = "Garnlav"
$ws.Range("G27").Value = "Alectoria sarmentosa"
$ws.Range("H27").Value = "(Ach.) Ach."
$ws.Range("Q27").Value = 451172.0902361136
$ws.Range("R27").Value = 7086726.569319103

# Row 28 (now holds the "Trådticka" / Climacocystis borealis record)
$ws.Range("A28").Value = 111936789
$ws.Range("B28").Value = 90087
$ws.Range("D28").Value = "LC"
$ws.Range("E28").Value = 3298
$ws.Range("F28").Value = "Trådticka"
$ws.Range("G28").Value = "Climacocystis borealis"
$ws.Range("H28").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q28").Value = 450955.1314140605
$ws.Range("R28").Value = 7087063.751596102

# Row 29 (now holds the "Granticka" / Porodaedalea chrysoloma record)
$ws.Range("A29").Value = 111936864
$ws.Range("B29").Value = 89423
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 5432
$ws.Range("F29").Value = "Granticka"
$ws.Range("G29").Value = "Porodaedalea chrysoloma"
$ws.Range("H29").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q29").Value = 451094.1385684713
$ws.Range("R29").Value = 7087212.607717684

# Row 30 (now holds the "Spillkråka" / Dryocopus martius record)
$ws.Range("A30").Value = 111936854
$ws.Range("B30").Value = 56414
$ws.Range("E30").Value = 100049
$ws.Range("F30").Value = "Spillkråka"
$ws.Range("G30").Value = "Dryocopus martius"
$ws.Range("K30:N30").ClearContents()
$ws.Range("Q30").Value = 450998.3386916541
$ws.Range("R30").Value = 7087288.958247212
$ws.Range("AC30").ClearContents()
